$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell without Excel coercing numeric-looking
# strings (e.g. "592.27") into floating point numbers, and without leaving a
# residual "Text" number-format style behind on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.386.83'
$ws.Range('E2').Value = '  +0.30%  '

Set-TextValue 'D3' '2.553.73'
$ws.Range('E3').Value = '  -2.39%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue 'D5' '592.27'
$ws.Range('E5').Value = '  +0.34%  '

Set-TextValue 'D6' '173.82'
$ws.Range('E6').Value = '  +4.82%  '

$ws.Range('E8').Value = '  -0.08%  '

Set-TextValue 'D9' '2.552.46'
$ws.Range('E9').Value = '  -2.39%  '

$ws.Range('E10').Value = '  +1.03%  '

Set-TextValue 'D11' '0.163'
$ws.Range('E11').Value = '  +1.40%  '

$ws.Range('E12').Value = '  -0.52%  '

Set-TextValue 'D13' '0.353'
$ws.Range('E13').Value = '  -3.02%  '

Set-TextValue 'D14' '27.14'
$ws.Range('E14').Value = '  -0.67%  '

Set-TextValue 'D15' '3.004.77'
$ws.Range('E15').Value = '  -2.74%  '

$ws.Range('E16').Value = '  -0.47%  '

Set-TextValue 'D17' '67.261.73'
$ws.Range('E17').Value = '  +0.16%  '

Set-TextValue 'D18' '2.550.96'
$ws.Range('E18').Value = '  -2.38%  '

Set-TextValue 'D19' '8.07'
$ws.Range('E19').Value = '  +3.54%  '

$ws.Range('E20').Value = '  -2.78%  '

Set-TextValue 'D21' '357.23'
$ws.Range('E21').Value = '  +0.61%  '

$ws.Range('E22').Value = '  -1.10%  '

Set-TextValue 'D23' '4.70'
$ws.Range('E23').Value = '  +1.21%  '

Set-TextValue 'D24' '2.01'
$ws.Range('E24').Value = '  +4.76%  '

$ws.Range('E25').Value = '  -0.03%  '

Set-TextValue 'D26' '70.17'
$ws.Range('E26').Value = '  +1.48%  '

Set-TextValue 'D27' '10.08'
$ws.Range('E27').Value = '  -3.83%  '

Set-TextValue 'D28' '2.685.93'
$ws.Range('E28').Value = '  -2.39%  '

Set-TextValue 'D29' '0.998'
$ws.Range('E29').Value = '  -0.07%  '

Set-TextValue 'D30' '0.0000100'
$ws.Range('E30').Value = '  +0.41%  '

Set-TextValue 'D31' '537.07'
$ws.Range('E31').Value = '  -1.06%  '

Set-TextValue 'D32' '8.26'
$ws.Range('E32').Value = '  +5.02%  '

Set-TextValue 'D33' '1.36'
$ws.Range('E33').Value = '  +1.62%  '

Set-TextValue 'D34' '1.87'
$ws.Range('E34').Value = '  -0.33%  '

Set-TextValue 'D35' '0.133'
$ws.Range('E35').Value = '  -0.26%  '

$ws.Range('E37').Value = '  -0.03%  '

Set-TextValue 'D38' '157.24'
$ws.Range('E38').Value = '  -1.18%  '

Set-TextValue 'D39' '18.84'
$ws.Range('E39').Value = '  -0.34%  '

Set-TextValue 'D40' '18.47'
$ws.Range('E40').Value = '  +1.22%  '

Set-TextValue 'D41' '0.358'
$ws.Range('E41').Value = '  -1.72%  '

$ws.Range('E42').Value = '  +0.55%  '

Set-TextValue 'D43' '5.22'
$ws.Range('E43').Value = '  +1.47%  '

Set-TextValue 'D44' '2.57'
$ws.Range('E44').Value = '  +6.56%  '

$ws.Range('E45').Value = '  +0.04%  '

$ws.Range('E46').Value = '  -0.86%  '

Set-TextValue 'D47' '151.54'
$ws.Range('E47').Value = '  +0.03%  '

Set-TextValue 'D50' '3.74'
$ws.Range('E50').Value = '  -0.72%  '

$ws.Range('E51').Value = '  +1.43%  '

# Rows 48 and 49: BabyDogeCoin and ARBITRUM swap rank position
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D48' '0.0₆0285'
$ws.Range('E48').Value = '  -3.90%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D49' '0.567'
$ws.Range('E49').Value = '  -1.46%  '
